# Experiment3-Tables.xlsx edit: rename smape_mean/smape_std/2-imputations-combined
# labels to their space-separated variants on "Sheet1", and update the last
# selected cell on "Sheet1" and "missing-1or1".

$wb = $excel.ActiveWorkbook

# --- Sheet1: update header / row labels -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C2").Value = "smape mean"
$ws1.Range("D2").Value = "smape std"
$ws1.Range("E2").Value = "smape mean"
$ws1.Range("F2").Value = "smape std"
$ws1.Range("G2").Value = "smape mean"
$ws1.Range("H2").Value = "smape std"

$ws1.Range("B3").Value = "2 imputations combined"
$ws1.Range("B9").Value = "2 imputations combined"
$ws1.Range("B15").Value = "2 imputations combined"

# --- missing-1or1: move the last saved selection -----------------------------------------
$ws4 = $wb.Worksheets.Item("missing-1or1")
$ws4.Activate()
$ws4.Range("H25").Select()

# --- Sheet1 stays the active/tab-selected sheet, with the last selection at L13 ----------
$ws1.Activate()
$ws1.Range("L13").Select()
